$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.842.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.99%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.249.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.95%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "395.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.28%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.580"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.21%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.247.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.94%  "

# Row 9
$ws.Range("E9").Value = "  -0.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.627"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "39.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.60%  "

# Row 12
$ws.Range("E12").Value = "  +9.70%  "

# Row 13
$ws.Range("E13").Value = "  +2.08%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.765.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.34%  "

# Row 15
$ws.Range("E15").Value = "  +3.34%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.36%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.252.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.95%  "

# Row 18
$ws.Range("E18").Value = "  -3.47%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.36%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "56.699.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.95%  "

# Row 21
$ws.Range("E21").Value = "  +1.10%  "

# Row 22
$ws.Range("E22").Value = "  +9.23%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "295.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.14%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.25"
$ws.Range("D25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "  -3.47%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.21%  "

# Row 28
$ws.Range("E28").Value = "  +0.66%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.91%  "

# Row 31
$ws.Range("E31").Value = "  -0.96%  "

# Row 32
$ws.Range("E32").Value = "  +0.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.74%  "

# Row 34
$ws.Range("E34").Value = "  -3.73%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "39.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.84%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0486"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.66%  "

# Row 37
$ws.Range("E37").Value = "  +2.04%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.11%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.61%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.03%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "137.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.35%  "

# Row 43
$ws.Range("E43").Value = "  +3.70%  "

# Row 44
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.81%  "

# Row 45
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.74%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.12%  "

# Row 47
$ws.Range("E47").Value = "  -3.36%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.13%  "

# Row 49
$ws.Range("E49").Value = "  +3.60%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.157.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.15%  "

# Row 51
$ws.Range("E51").Value = "  -5.51%  "
